$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (Wins / Losses / Ties) in AD1:AF1.
# Set the values first, then copy the existing header formatting
# (bold font, thin border, centered alignment) from AC1 onto them so
# the new cells share the same style as the rest of row 1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Season record values for every data row (2-47): 107 wins, 55 losses, 0 ties.
$lastRow = 47
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 107
    $ws.Cells.Item($r, 31).Value = 55
    $ws.Cells.Item($r, 32).Value = 0
}
